{"js": "// Revise the sentence describing the web app / training app development,\n// per commit \"Revise app and thesis also add load method\".\n//\n// Before (concatenated run text):\n//   \" Aplikasi untuk mengumpulkan dan menirukan suara akan dikembangkan dalam website\"\n//   \" dan \"\n//   \"aplikasi untuk melatih dikembangkan dalam \"\n//   [bookmarkStart/_GoBack][bookmarkEnd]\n//   \"command promt\"\n//   \".\"\n//\n// After:\n//   \" Aplikasi untuk mengumpulkan dan menirukan suara dikembangkan\"\n//   [bookmarkStart/_GoBack][bookmarkEnd]\n//   \" berbasis web\"\n//   \" dan \"\n//   \"aplikasi untuk melatih dikembangkan\"\n//   \" berbasis teks\"\n//   \".\"\n\nconst body = context.document.body;\n\n// 1) \"...suara akan dikembangkan dalam website\" -> \"...suara dikembangkan berbasis web\"\nlet res1 = body.search(\"akan dikembangkan dalam website\", { matchCase: true });\nres1.load(\"items\");\nawait context.sync();\nif (res1.items.length > 0) {\n  res1.items[0].insertText(\"dikembangkan berbasis web\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"...melatih dikembangkan dalam \" -> \"...melatih dikembangkan\" (drop the trailing \"dalam \")\nlet res2 = body.search(\"melatih dikembangkan dalam \", { matchCase: true });\nres2.load(\"items\");\nawait context.sync();\nif (res2.items.length > 0) {\n  res2.items[0].insertText(\"melatih dikembangkan\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) \"command promt\" -> \" berbasis teks\" (leading space replaces the one dropped in step 2)\nlet res3 = body.search(\"command promt\", { matchCase: true });\nres3.load(\"items\");\nawait context.sync();\nif (res3.items.length > 0) {\n  res3.items[0].insertText(\" berbasis teks\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4) The _GoBack bookmark marks the author's last edit point, which in the revised text sits\n//    right after \"...suara dikembangkan\" (before \" berbasis web\"). Re-anchor it there.\nconst oldGoBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\noldGoBack.load(\"isNullObject\");\nawait context.sync();\nif (!oldGoBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\nlet res4 = body.search(\"suara dikembangkan\", { matchCase: true });\nres4.load(\"items\");\nawait context.sync();\nif (res4.items.length > 0) {\n  res4.items[0].getRange(Word.RangeLocation.end).insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Revise the sentence describing the web app / training app development,\n# per commit \"Revise app and thesis also add load method\".\n#\n# Before:\n#   \" Aplikasi untuk mengumpulkan dan menirukan suara akan dikembangkan dalam website\"\n#   \" dan \"\n#   \"aplikasi untuk melatih dikembangkan dalam \"\n#   [bookmarkStart/_GoBack][bookmarkEnd]\n#   \"command promt\"\n#   \".\"\n#\n# After:\n#   \" Aplikasi untuk mengumpulkan dan menirukan suara dikembangkan\"\n#   [bookmarkStart/_GoBack][bookmarkEnd]\n#   \" berbasis web\"\n#   \" dan \"\n#   \"aplikasi untuk melatih dikembangkan\"\n#   \" berbasis teks\"\n#   \".\"\n\n$d = $word.ActiveDocument\n\n# 1) \"...suara akan dikembangkan dalam website\" -> \"...suara dikembangkan berbasis web\"\n$r1 = $d.Content\n$r1.Find.Execute(\"akan dikembangkan dalam website\", $false, $false, $false, $false, $false, $true, 1, $false, \"dikembangkan berbasis web\", 2)\n\n# 2) \"...melatih dikembangkan dalam \" -> \"...melatih dikembangkan\" (drop the trailing \"dalam \")\n$r2 = $d.Content\n$r2.Find.Execute(\"melatih dikembangkan dalam \", $false, $false, $false, $false, $false, $true, 1, $false, \"melatih dikembangkan\", 2)\n\n# 3) \"command promt\" -> \" berbasis teks\" (leading space restores the one dropped in step 2)\n$r3 = $d.Content\n$r3.Find.Execute(\"command promt\", $false, $false, $false, $false, $false, $true, 1, $false, \" berbasis teks\", 2)\n\n# 4) The _GoBack bookmark marks the point of the author's last edit, which in the revised\n#    text sits right after \"... suara dikembangkan\" (before \" berbasis web\"). Re-anchor it there.\n$r4 = $d.Content\n$r4.Find.Execute(\"suara dikembangkan\")\nif ($r4.Find.Found) {\n    $bmRange = $d.Range($r4.End, $r4.End)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
